$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# --- Constants sheet: remove the blank row 2, which shifts every row below it up by one ---
$wsConstants.Rows.Item(2).Delete()

# --- Recreate the navigation / selection state captured in the saved workbook ---
# Visit Assets, then Constants, then Settings last, leaving Settings as the
# active sheet (and active cell A6) when the workbook is saved, matching the
# activeTab / tabSelected / selection changes in the diff.
$wsAssets.Activate()
$wsAssets.Range("A2").Select()

$wsConstants.Activate()
$wsConstants.Range("A12").Select()

$wsSettings.Activate()
$wsSettings.Range("A6").Select()
